# Apply "repull data, push all data, mean calculation" update:
# Column F ("dSF") gets refreshed values for most rows (a recomputed
# differential that is no longer simply mirroring column E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = -5
    4  = 1
    5  = 2
    6  = -3
    8  = 7
    11 = -3
    12 = -6
    13 = 3
    14 = -1
    15 = 6
    16 = -1
    17 = -6
    18 = 1
    19 = 4
    20 = -1
    21 = -2
    22 = 8
    23 = 3
    24 = 3
    25 = -2
    26 = 5
    27 = 3
    28 = 6
    30 = 3
    31 = -1
    32 = -5
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
